# Apply updated data-masking logic to the Student Information sheet.
# Masks column A (Name), column E (Phone Number) and column G (Email)
# for rows 2-31, keeping the first and last character visible and
# replacing everything in between with asterisks. For email addresses,
# only the local part (before the "@") is masked; the domain is left
# untouched.

function Mask-Value([string]$value) {
    if ($value.Length -le 2) {
        return $value
    }
    $first = $value.Substring(0, 1)
    $last = $value.Substring($value.Length - 1, 1)
    $middle = "".PadLeft($value.Length - 2, '*')
    return "$first$middle$last"
}

function Mask-Email([string]$value) {
    $atIndex = $value.IndexOf('@')
    if ($atIndex -lt 0) {
        return Mask-Value $value
    }
    $local = $value.Substring(0, $atIndex)
    $domain = $value.Substring($atIndex)
    return (Mask-Value $local) + $domain
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $nameCell = $ws.Cells.Item($row, 1)   # Column A - Name
    $phoneCell = $ws.Cells.Item($row, 5)  # Column E - Phone Number
    $emailCell = $ws.Cells.Item($row, 7)  # Column G - Email

    $nameCell.Value = Mask-Value $nameCell.Value()
    $phoneCell.Value = Mask-Value $phoneCell.Value()
    $emailCell.Value = Mask-Email $emailCell.Value()
}
